# mod port og ch8, add paper SPY/TLT
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Portfolio holdings table: update TSLA (row 24) quantity from 16 to 14
$ws.Range("C24").Value = 14

# Move the active selection to C25 (last cell of the table)
$ws.Range("C25").Select()
